$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'23.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.366"
$ws.Range("D4").Style = "Normal"
$ws.Range("D6").Value = "'3.377"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.486"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Value = "'0.9197"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1421"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07408"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03096"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03068"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09355"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.875"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001556"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04706"
$ws.Range("D17").Style = "Normal"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005940"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.005961"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001242"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.004719"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.00008804"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21NitroExNTXBestin24h"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.595"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.158"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("D25").Value = "'0.3232"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1329"
$ws.Range("D26").Style = "Normal"
$ws.Range("D40").Value = "'0.03867"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1068"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002781"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003113"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("D44").Value = "'0.008568"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005246"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.7100"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.001737"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "47BOLOBOLO"
